$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the c27a0754 (handback) row as handed back ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status + handback datetime for the c27a0754 row ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-01-13 08:32:43"

# --- de-de sheet: status + handback datetime for the c27a0754 row ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-01-13 08:33:18"
